## Applies the "i am a student" edit:
##  1. Splits the existing "Palak git scm" run into two runs ("Palak git "
##     and "scm"), wrapping the word "scm" in spell-check proofing marks
##     (w:proofErr spellStart/spellEnd), as Word does for words it does
##     not recognise.
##  2. Adds a new paragraph "Pappaaaa mummma", with "Pappaaaa" and
##     "mummma" each wrapped in their own spell-check proofing marks.
##  3. Adds a new, completely empty paragraph after that.

$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Paragraph 1: "Palak git scm" -> "Palak git " + proofed "scm" -----------
$para1 = $d.Paragraphs(1).Range
$xmlPara1 = "<w:p $wns>" +
              "<w:r><w:t xml:space='preserve'>Palak git </w:t></w:r>" +
              "<w:proofErr w:type='spellStart'/>" +
              "<w:r><w:t>scm</w:t></w:r>" +
              "<w:proofErr w:type='spellEnd'/>" +
            "</w:p>"
$para1.InsertXML($xmlPara1)

# --- Append two new paragraphs right after paragraph 1 -----------------------
$tail = $d.Paragraphs(1).Range
$tail.InsertParagraphAfter()
$tail.InsertParagraphAfter()

# --- Paragraph 2: "Pappaaaa mummma" with proofing marks on each word --------
$para2 = $d.Paragraphs(2).Range
$xmlPara2 = "<w:p $wns>" +
              "<w:proofErr w:type='spellStart'/>" +
              "<w:r><w:t>Pappaaaa</w:t></w:r>" +
              "<w:proofErr w:type='spellEnd'/>" +
              "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
              "<w:proofErr w:type='spellStart'/>" +
              "<w:r><w:t>mummma</w:t></w:r>" +
              "<w:proofErr w:type='spellEnd'/>" +
            "</w:p>"
$para2.InsertXML($xmlPara2)

# --- Paragraph 3: a new, fully empty paragraph -------------------------------
$para3 = $d.Paragraphs(3).Range
$para3.InsertXML("<w:p $wns/>")
